# Update column C ("Förändrad") date values from 45175 to 45177
# for every data row (rows 2 through 123) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 123; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45177
    }
}
